# Handles float input without breaking stuff
#
# Updates the computed marksheet summary (rows 10-12), recolors the
# "Student Ans" column of the answer grid to reflect correct / incorrect /
# not-attempted answers, fills in the second question block's "Student Ans"
# column (D16:D18), and drops the now-unused third question block
# (columns G:H) together with the now-empty D/E rows below row 18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Summary table (rows 10-12)
# ---------------------------------------------------------------------

# Give the row labels in column A the same bold/bordered style used by
# the header row above them (row 9) / E11.
$labelStyleSrc = $ws.Range("A9")
foreach ($r in 10, 11, 12) {
    $labelStyleSrc.Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

$ws.Range("B10").Value = 22
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 88
$ws.Range("C12").Value = -3
$ws.Range("E12").Value = "85/112"

# ---------------------------------------------------------------------
# 2. Answer grid - color / fill in "Student Ans" columns
# ---------------------------------------------------------------------

# Style donors already present in the sheet.
$greenSrc = $ws.Range("B10")   # "correct" (green) style
$redSrc   = $ws.Range("C10")   # "incorrect" (red) style

# Block 1 : columns A (student) / B (correct), rows 16-40
# value, outcome ("green"/"red"/"none")
$block1 = @{
    16 = @("Option A", "green")
    17 = @("Option D", "green")
    18 = @("Option B", "green")
    19 = @("Option C", "green")
    20 = @($null,      "none")
    21 = @("Option C", "green")
    22 = @("Option D", "green")
    23 = @("Option D", "green")
    24 = @($null,      "none")
    25 = @("Option A", "green")
    26 = @("Option C", "green")
    27 = @("Option A", "green")
    28 = @("Option D", "green")
    29 = @("Option D", "green")
    30 = @("Option B", "green")
    31 = @("Option D", "green")
    32 = @("Option C", "green")
    33 = @("Option D", "green")
    34 = @("Option A", "red")
    35 = @("Option D", "green")
    36 = @("Option D", "red")
    37 = @("Option D", "red")
    38 = @("Option A", "green")
    39 = @("Option D", "green")
    40 = @($null,      "none")
}

foreach ($r in 16..40) {
    $info = $block1[$r]
    $value = $info[0]
    $outcome = $info[1]

    if ($outcome -eq "green") {
        $greenSrc.Copy() | Out-Null
        $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
        $ws.Range("A$r").Value = $value
    } elseif ($outcome -eq "red") {
        $redSrc.Copy() | Out-Null
        $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
        $ws.Range("A$r").Value = $value
    }
    # "none" rows (20, 24, 40) are left exactly as-is: empty, unattempted.
}

# Block 2 : columns D (student) / E (correct) - only rows 16-18 still
# have data; the student got every one of them right.
foreach ($r in 16, 17, 18) {
    $greenSrc.Copy() | Out-Null
    $ws.Range("D$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("D$r").Value = $ws.Range("E$r").Value2
}

# ---------------------------------------------------------------------
# 3. Remove the leftover/unused cells:
#    - D/E below row 18 (block 2 only had 3 questions)
#    - the whole third question block in columns G/H
# ---------------------------------------------------------------------

$ws.Range("D19:E40").Clear() | Out-Null
$ws.Range("G15:H40").Clear() | Out-Null

Write-Host "edit applied"
